# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullet list to use
# impact-focused accomplishment statements, trimming from 6 bullets
# down to 4.

$d = $word.ActiveDocument

# Locate the "Impact" sub-heading that lives under the
# "KEY ACHIEVEMENTS AND IMPACT" section so edits are scoped to that
# block only (some of the bullet text also appears earlier in the
# document, under Professional Experience, and must be left alone).
$scopeStart = $d.Content
$found = $scopeStart.Find.Execute("KEY ACHIEVEMENTS AND IMPACT", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$sectionRange = $d.Range($scopeStart.End, $d.Content.End)

# Replace the first three bullets in place (1:1 text swap).
$r1 = $sectionRange.Duplicate
$r1.Find.Execute(
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Revenue generation: Delivered `$4.9M additional revenue through optimization",
    2)

$r2 = $sectionRange.Duplicate
$r2.Find.Execute(
    "• Delivered `$4.9M additional revenue through continuous testing and optimization, increased conversion rates by 23%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• 23% conversion rate improvement",
    2)

$r3 = $sectionRange.Duplicate
$r3.Find.Execute(
    "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis",
    2)

# Replace the 6th bullet's text before removing the two bullets
# (#4 and #5) that sit between it and bullet #3.
$r6 = $sectionRange.Duplicate
$r6.Find.Execute(
    "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations",
    2)

# Now delete the two whole paragraphs that fall between the
# "Executive authority" bullet and the "Platform impact" bullet:
#   • Developed longitudinal data analysis methods ...
#   • Discovered systematic race coding errors ...
foreach ($needle in @(
        "• Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality",
        "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%"
    )) {
    $rd = $d.Range($scopeStart.End, $d.Content.End).Duplicate
    $rd.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($rd.Find.Found) {
        $para = $rd.Paragraphs(1).Range
        $para.Delete()
    }
}
